$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 gets a new date (A3) matching the date-style used by A2, and a
# work-hours note (B3) referencing a new shared string.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 43980
$ws.Range("B3").Value = "4 timer 33 minutter"

# Update the current selection / cursor position and zoom level to match
# the author's final view state.
$ws.Range("D7").Select()
$excel.ActiveWindow.Zoom = 160
